$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in column F, matching the style of the existing headers (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats

$timestamps = @(
    "2021-10-05 13:40:52.029505",
    "2021-10-05 13:40:52.029518",
    "2021-10-05 13:40:52.029521",
    "2021-10-05 13:40:52.029524",
    "2021-10-05 13:40:52.029526",
    "2021-10-05 13:40:52.029529",
    "2021-10-05 13:40:52.029531",
    "2021-10-05 13:40:52.029534",
    "2021-10-05 13:40:52.029536",
    "2021-10-05 13:40:52.029539",
    "2021-10-05 13:40:52.029541",
    "2021-10-05 13:40:52.029543",
    "2021-10-05 13:40:52.029546",
    "2021-10-05 13:40:52.029549",
    "2021-10-05 13:40:52.029551",
    "2021-10-05 13:40:52.029553",
    "2021-10-05 13:40:52.029556",
    "2021-10-05 13:40:52.029559",
    "2021-10-05 13:40:52.029561"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
